$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2024-05-27", "11:48:23", "-", "Etiquetadora", "-", "-", "-", "11:48:26", "0:00:03"),
    @("2024-05-27", "11:48:27", "-", "No coloca bien el sealling", "-", "-", "-", "11:48:29", "0:00:02"),
    @("2024-05-27", "12:13:45", "-", "-", "-", "Robot no coloca bien filter en palet", "-", "12:13:47", "0:00:02"),
    @("2024-05-27", "12:14:07", "-", "-", "-", "NOK Soldadura metal", "-", "12:14:08", "0:00:01"),
    @("2024-05-27", "12:14:24", "-", "-", "-", "Robot no coloca bien filter en palet", "-", "12:14:28", "0:00:04"),
    @("2024-05-27", "12:14:26", "-", "-", "-", "Robot no coloca bien filter en palet", "-", "12:14:28", "0:00:02"),
    @("2024-05-27", "12:14:33", "-", "-", "-", "Robot no coloca bien filter en palet", "-", "12:14:35", "0:00:02"),
    @("2024-05-27", "12:46:23", "-", "Cámara no detecta Pcb", "-", "-", "-", "12:46:25", "0:00:02"),
    @("2024-05-27", "12:46:28", "-", "Cámara no detecta busbar", "-", "-", "-", "12:46:30", "0:00:02")
)

$startRow = 300
$endRow = $startRow + $data.Count - 1

# Column A holds values that look like ISO dates ("2024-05-27"). Excel's
# auto-detection would otherwise silently convert these into date serial
# numbers. Temporarily mark the column as Text so the values are stored
# verbatim as strings, then clear the formatting again afterwards so the
# cells are left with the default (General) style, matching the rest of
# the sheet which carries no explicit formatting.
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}

$dateRange.ClearFormats()
